$wb = $excel.ActiveWorkbook

# --- Sheet "2023" (sheet1): add a new timeline entry, resize col A, set selection ---
$ws2023 = $wb.Worksheets.Item("2023")
$ws2023.Range("A2").Value = "Started working on LLMs at Barclays"
$ws2023.Range("B2").Value = "October"
$ws2023.Columns.Item(1).ColumnWidth = 30
$ws2023.Range("C5").Select()

# --- Sheet "2020" (sheet4): it was the active tab before; now it just loses focus
# and gets a different remembered selection. Set its selection while it's still
# active so it doesn't flip back to being the active tab afterwards. ---
$ws2020 = $wb.Worksheets.Item("2020")
$ws2020.Range("B8").Select()

# --- Sheet "2021" (sheet3): add a new timeline entry, resize col A, becomes the
# active tab, with its own remembered selection. Do this last so it ends up
# being the active sheet in the saved workbook. ---
$ws2021 = $wb.Worksheets.Item("2021")
$ws2021.Range("A2").Value = "Did Internship at RTG Computational Cognition, Germany"
$ws2021.Range("B2").Value = "June"
$ws2021.Columns.Item(1).ColumnWidth = 47.25
$ws2021.Activate()
$ws2021.Range("A12").Select()
